$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, pushing the existing rows 59:83 down to 60:84
# (this also extends the used range / dimension to row 84, matching row 83's old
# style carrying down correctly for the date column).
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new weekly record.
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = "Macroferia Regional de Talca"
$ws.Range("C59").Value = "Maule"
$ws.Range("D59").Value = 44876
$ws.Range("E59").Value = 7
$ws.Range("F59").Value = 300000000
$ws.Range("G59").Value = "Espárragos"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 3000
$ws.Range("K59").Value = 1000
$ws.Range("L59").Value = 1000
$ws.Range("M59").Value = 1000
$ws.Range("N59").Value = '$/kilo'
$ws.Range("O59").Value = "Provincia de Linares"
$ws.Range("P59").Value = 1000
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"
